# Apply updated crypto price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.703.94"
$ws.Range("E2").Value = "  -1.82%  "
$ws.Range("D3").Value = "1.754.54"
$ws.Range("E3").Value = "  -2.51%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "'324.75"
$ws.Range("E5").Value = "  -4.12%  "
$ws.Range("D6").Value = "'0.9999"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").Value = "'0.4481"
$ws.Range("E7").Value = "  -2.29%  "
$ws.Range("D8").Value = "'0.3689"
$ws.Range("E8").Value = "  -2.27%  "
$ws.Range("D9").Value = "'45.25"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").Value = "'0.07488"
$ws.Range("E10").Value = "  -1.84%  "
$ws.Range("D11").Value = "'1.122"
$ws.Range("E11").Value = "  -2.04%  "
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("D13").Value = "'21.70"
$ws.Range("E13").Value = "  -3.06%  "
$ws.Range("E14").Value = "  -2.59%  "
$ws.Range("D15").Value = "'7.276"
$ws.Range("E15").Value = "  -3.15%  "
$ws.Range("D16").Value = "1.750.72"
$ws.Range("E16").Value = "  -2.81%  "
$ws.Range("D17").Value = "'0.00001072"
$ws.Range("E17").Value = "  -1.65%  "
$ws.Range("D18").Value = "'88.22"
$ws.Range("E18").Value = "  +8.34%  "
$ws.Range("D19").Value = "'0.06220"
$ws.Range("E19").Value = "  -7.65%  "
$ws.Range("D20").Value = "'0.9999"
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("D21").Value = "'17.15"
$ws.Range("E21").Value = "  -1.48%  "
$ws.Range("D22").Value = "'6.168"
$ws.Range("E22").Value = "  -3.72%  "
$ws.Range("D23").Value = "'0.5293"
$ws.Range("E23").Value = "  -3.83%  "
$ws.Range("D24").Value = "27.729.16"
$ws.Range("E24").Value = "  -1.76%  "
$ws.Range("D25").Value = "'11.64"
$ws.Range("E25").Value = "  -1.88%  "
$ws.Range("D26").Value = "'2.325"
$ws.Range("E26").Value = "  -4.08%  "
$ws.Range("D27").Value = "'20.62"
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("D28").Value = "'153.95"
$ws.Range("E28").Value = "  +1.13%  "
$ws.Range("D29").Value = "'2.357"
$ws.Range("E29").Value = "  -0.28%  "
$ws.Range("D30").Value = "1.949.38"
$ws.Range("E30").Value = "  -2.86%  "
$ws.Range("D31").Value = "'127.86"
$ws.Range("E31").Value = "  -3.82%  "
$ws.Range("D32").Value = "'1.219"
$ws.Range("E32").Value = "  -1.23%  "
$ws.Range("D33").Value = "'5.737"
$ws.Range("E33").Value = "  -1.56%  "
$ws.Range("D34").Value = "'0.09273"
$ws.Range("E34").Value = "  -2.31%  "
$ws.Range("D35").Value = "'3.632"
$ws.Range("E35").Value = "  -10.01%  "
$ws.Range("D36").Value = "'12.67"
$ws.Range("E36").Value = "  +4.95%  "
$ws.Range("B37").Value = "Algorand"
$ws.Range("C37").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D37").Value = "'0.2169"
$ws.Range("E37").Value = "  -5.59%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.02313"
$ws.Range("E38").Value = "  -1.50%  "
$ws.Range("D39").Value = "'5.098"
$ws.Range("E39").Value = "  -2.99%  "
$ws.Range("D40").Value = "'0.6487"
$ws.Range("E40").Value = "  -1.68%  "
$ws.Range("D41").Value = "'0.06120"
$ws.Range("E41").Value = "  -3.67%  "
$ws.Range("D42").Value = "'1.201"
$ws.Range("E42").Value = "  -2.69%  "
$ws.Range("D43").Value = "'7.967"
$ws.Range("E43").Value = "  -5.25%  "
$ws.Range("E44").Value = "  -4.41%  "
$ws.Range("D45").Value = "'0.9996"
$ws.Range("E45").Value = "  -0.23%  "
$ws.Range("D46").Value = "'13.81"
$ws.Range("E46").Value = "  -2.39%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "'3.755"
$ws.Range("E47").Value = "  -2.92%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Value = "'0.5949"
$ws.Range("E48").Value = "  -2.56%  "
$ws.Range("D49").Value = "'126.25"
$ws.Range("E49").Value = "  -3.47%  "
$ws.Range("D50").Value = "'1.979"
$ws.Range("E50").Value = "  -2.49%  "
$ws.Range("D51").Value = "'0.06901"
$ws.Range("E51").Value = "  -3.42%  "
